$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Save" header in H1, matching the style used by the other headers (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Values for the new "Save" column, one per data row (rows 2-17)
$saveValues = @(0, 0, 1, 0, 0, 0, 1, 1, 0, 0, 0, 0, 0, 1, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
